$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 45939
$ws.Cells.Item(2, 3).Value = 0.3413142033632344
$ws.Cells.Item(3, 1).Value = 45940
$ws.Cells.Item(3, 3).Value = 0.3190980429478851
$ws.Cells.Item(4, 1).Value = 45943
$ws.Cells.Item(4, 3).Value = 0.3302904728403517
$ws.Cells.Item(5, 1).Value = 45944
$ws.Cells.Item(5, 3).Value = 0.3539499120406343
$ws.Cells.Item(6, 1).Value = 45945
$ws.Cells.Item(6, 3).Value = 0.3678546466009929
$ws.Cells.Item(7, 1).Value = 45946
$ws.Cells.Item(7, 3).Value = 0.3793855730025166
$ws.Cells.Item(8, 1).Value = 45947
$ws.Cells.Item(8, 3).Value = 0.3635105148825148
$ws.Cells.Item(9, 1).Value = 45950
$ws.Cells.Item(9, 3).Value = 0.3910460700388858
$ws.Cells.Item(10, 1).Value = 45951
$ws.Cells.Item(10, 3).Value = 0.3691685531396178
$ws.Cells.Item(11, 1).Value = 45952
$ws.Cells.Item(11, 3).Value = 0.3572080456248493
$ws.Cells.Item(12, 1).Value = 45953
$ws.Cells.Item(12, 3).Value = 0.319388135122457
$ws.Cells.Item(13, 1).Value = 45954
$ws.Cells.Item(13, 3).Value = 0.2778647663256245
$ws.Cells.Item(14, 1).Value = 45957
$ws.Cells.Item(14, 3).Value = 0.2828287453755624
$ws.Cells.Item(15, 1).Value = 45958
$ws.Cells.Item(15, 3).Value = 0.2901432403350007
$ws.Cells.Item(16, 1).Value = 45959
$ws.Cells.Item(16, 3).Value = 0.2736733886462862
$ws.Cells.Item(17, 1).Value = 45960
$ws.Cells.Item(17, 3).Value = 0.2771206740010897
$ws.Cells.Item(18, 1).Value = 45961
$ws.Cells.Item(18, 3).Value = 0.2949025682274686
$ws.Cells.Item(19, 1).Value = 45964
$ws.Cells.Item(19, 3).Value = 0.3267378078632872
$ws.Cells.Item(20, 1).Value = 45965
$ws.Cells.Item(20, 3).Value = 0.3529978613329278
$ws.Cells.Item(21, 1).Value = 45966
$ws.Cells.Item(21, 3).Value = 0.3574175390662176
$ws.Cells.Item(22, 1).Value = 45967
$ws.Cells.Item(22, 3).Value = 0.4091398488479487
$ws.Cells.Item(23, 1).Value = 45968
$ws.Cells.Item(23, 3).Value = 0.4172458997115651
$ws.Cells.Item(24, 1).Value = 45971
$ws.Cells.Item(24, 3).Value = 0.4188715865954915
$ws.Cells.Item(25, 1).Value = 45972
$ws.Cells.Item(25, 3).Value = 0.4092736853072489
$ws.Cells.Item(26, 1).Value = 45973
$ws.Cells.Item(26, 3).Value = 0.4192898713098626
$ws.Cells.Item(27, 1).Value = 45974
$ws.Cells.Item(27, 3).Value = 0.4074392460332971
$ws.Cells.Item(28, 1).Value = 45975
$ws.Cells.Item(28, 3).Value = 0.4118779145227532
$ws.Cells.Item(29, 1).Value = 45978
$ws.Cells.Item(29, 3).Value = 0.4154290698515039
$ws.Cells.Item(30, 1).Value = 45979
$ws.Cells.Item(30, 3).Value = 0.4171573363290887
$ws.Cells.Item(31, 1).Value = 45980
$ws.Cells.Item(31, 3).Value = 0.4171605251775842
